$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("2025-02-27 16:50", "무색 큐브 조각", 58)
    ,@("2025-02-27 16:50", "모순의 결정체", 61705)
    ,@("2025-02-27 16:50", "레어 소울 결정", 82777)
    ,@("2025-02-27 16:50", "유니크 소울 결정", 253937)
    ,@("2025-02-27 16:50", "레전더리 소울 결정", 1384216)
    ,@("2025-02-27 16:50", "에픽 소울 결정", 4369810)
    ,@("2025-02-27 16:50", "태초 소울 결정", 46744787)
    ,@("2025-02-28 23:26", "무색 큐브 조각", 62)
    ,@("2025-02-28 23:26", "모순의 결정체", 63535)
    ,@("2025-02-28 23:26", "레어 소울 결정", 80446)
    ,@("2025-02-28 23:26", "유니크 소울 결정", 291183)
    ,@("2025-02-28 23:26", "레전더리 소울 결정", 1306346)
    ,@("2025-02-28 23:26", "에픽 소울 결정", 4361094)
    ,@("2025-02-28 23:26", "태초 소울 결정", 44467109)
)

$startRow = 323
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

Write-Host ("New dimension: " + $ws.UsedRange.Address())
